$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-09-02 -> 2023-09-03, serial 45171 -> 45172) for every data row.
$ws.Range("C2:C203").Value = 45172
